$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "TEST"
$r = $ws.Range("A1")
$r.Value
